{"js": "// Replace the old \"Dates de la campanya ...\" sentence with the new wording\n// in every paragraph of the document body where it occurs (4 occurrences).\nconst oldText = \"Dates de la campanya Constel\\u00B7laci\\u00F3 de Leo 2022: 14-23 d'abril, 14-23 de maig\";\nconst newText = \"Dates de la campanya 2022 en qu\\u00E8 usem la constel\\u00B7laci\\u00F3, Constel\\u00B7laci\\u00F3 de Leo 14-23 d'abril, 14-23 de maig\";\n\nconst results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(newText, Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# Replace the old \"Dates de la campanya ...\" sentence with the new wording\n# everywhere it occurs in the document body (4 occurrences).\n$d = $word.ActiveDocument\n\n$oldText = \"Dates de la campanya Constel\" + [char]0x00B7 + \"laci\" + [char]0x00F3 + \" de Leo 2022: 14-23 d'abril, 14-23 de maig\"\n$newText = \"Dates de la campanya 2022 en qu\" + [char]0x00E8 + \" usem la constel\" + [char]0x00B7 + \"laci\" + [char]0x00F3 + \", Constel\" + [char]0x00B7 + \"laci\" + [char]0x00F3 + \" de Leo 14-23 d'abril, 14-23 de maig\"\n\n# Use Find.Execute to locate each occurrence, then assign straight to\n# Range.Text (instead of passing the replacement through Find.Execute's\n# ReplaceWith argument) so AutoCorrect's \"smart quotes\" substitution does\n# not mangle the straight apostrophes in the replacement text.\n$found = $true\n$count = 0\nwhile ($found -and $count -lt 20) {\n    $searchRange = $d.Content\n    $searchRange.Start = 0\n    $find = $searchRange.Find\n    $find.ClearFormatting()\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $found = $find.Execute($oldText)\n    if ($found) {\n        $searchRange.Text = $newText\n        $count = $count + 1\n    }\n}\n"}
